# Auto-generated: apply scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# row 32
$ws_ALC.Range("H32").Value = 9369.111000000001
$ws_ALC.Range("I32").Value = 9720
$ws_ALC.Range("J32").Value = 8667.333000000001
$ws_ALC.Range("K32").Value = 9720
$ws_ALC.Range("L32").Value = 8667.333000000001
$ws_ALC.Range("M32").Value = -9394
$ws_ALC.Range("N32").Value = -9319.333000000001
# row 129
$ws_ALC.Range("H129").Value = 2321.9092
$ws_ALC.Range("I129").Value = 1541.8334
$ws_ALC.Range("K129").Value = 4625.5002
$ws_ALC.Range("M129").Value = 374.4997999999996
# row 132
$ws_ALC.Range("H132").Value = 1085.1471
$ws_ALC.Range("I132").Value = 1059.8387
$ws_ALC.Range("J132").Value = 1346.6666
$ws_ALC.Range("K132").Value = 3179.5161
$ws_ALC.Range("L132").Value = 4039.9998
$ws_ALC.Range("M132").Value = -649.5160999999998
$ws_ALC.Range("N132").Value = -9099.9998
# row 135
$ws_ALC.Range("H135").Value = 947.28
$ws_ALC.Range("I135").Value = 687.85
$ws_ALC.Range("K135").Value = 6190.650000000001
$ws_ALC.Range("M135").Value = -3655.650000000001
# row 141
$ws_ALC.Range("H141").Value = 1018.125
$ws_ALC.Range("I141").Value = 1018.125
$ws_ALC.Range("J141").Value = 0
$ws_ALC.Range("K141").Value = 3054.375
$ws_ALC.Range("L141").Value = 0
$ws_ALC.Range("M141").Value = 2125.625
$ws_ALC.Range("N141").ClearContents()

$ws_ARM = $wb.Worksheets.Item("ARM")
# row 26
$ws_ARM.Range("H26").Value = 2833
$ws_ARM.Range("I26").Value = 1749.5
$ws_ARM.Range("J26").Value = 5000
$ws_ARM.Range("K26").Value = 1749.5
$ws_ARM.Range("L26").Value = 5000
$ws_ARM.Range("M26").Value = -1419.5
$ws_ARM.Range("N26").Value = -5660
# row 29
$ws_ARM.Range("H29").Value = 9999.5
$ws_ARM.Range("J29").Value = 9999.5
$ws_ARM.Range("L29").Value = 9999.5
$ws_ARM.Range("N29").Value = -10615.5
# row 61
$ws_ARM.Range("H61").Value = 3404.6667
$ws_ARM.Range("I61").Value = 3319.7705
$ws_ARM.Range("K61").Value = 3319.7705
$ws_ARM.Range("M61").Value = -3107.7705
# row 132
$ws_ARM.Range("H132").Value = 2001.9546
$ws_ARM.Range("I132").Value = 1286.7869
$ws_ARM.Range("K132").Value = 3860.3607
$ws_ARM.Range("M132").Value = -1330.3607
# row 136
$ws_ARM.Range("H136").Value = 3404.6667
$ws_ARM.Range("I136").Value = 3319.7705
$ws_ARM.Range("K136").Value = 9959.3115
$ws_ARM.Range("M136").Value = -7409.3115

$ws_BSM = $wb.Worksheets.Item("BSM")
# row 134
$ws_BSM.Range("H134").Value = 7002.6
$ws_BSM.Range("I134").Value = 1999.8334
$ws_BSM.Range("J134").Value = 14506.75
$ws_BSM.Range("K134").Value = 5999.5002
$ws_BSM.Range("L134").Value = 43520.25
$ws_BSM.Range("M134").Value = -3464.5002
$ws_BSM.Range("N134").Value = -48590.25

$ws_CRP = $wb.Worksheets.Item("CRP")
# row 22
$ws_CRP.Range("H22").Value = 4031.2
$ws_CRP.Range("J22").Value = 4031.2
$ws_CRP.Range("L22").Value = 4031.2
$ws_CRP.Range("N22").Value = -4731.2
# row 31
$ws_CRP.Range("H31").Value = 31800.139
$ws_CRP.Range("I31").Value = 1875.5385
$ws_CRP.Range("K31").Value = 1875.5385
$ws_CRP.Range("M31").Value = -1580.5385
# row 32
$ws_CRP.Range("H32").Value = 5502.5
$ws_CRP.Range("I32").Value = 2505
$ws_CRP.Range("J32").Value = 8500
$ws_CRP.Range("K32").Value = 2505
$ws_CRP.Range("L32").Value = 8500
$ws_CRP.Range("M32").Value = -2189
$ws_CRP.Range("N32").Value = -9132
# row 34
$ws_CRP.Range("H34").Value = 31800.139
$ws_CRP.Range("I34").Value = 1875.5385
$ws_CRP.Range("K34").Value = 1875.5385
$ws_CRP.Range("M34").Value = -1673.5385
# row 36
$ws_CRP.Range("H36").Value = 3024.5
$ws_CRP.Range("I36").Value = 50
$ws_CRP.Range("J36").Value = 5999
$ws_CRP.Range("K36").Value = 50
$ws_CRP.Range("L36").Value = 5999
$ws_CRP.Range("M36").Value = 338
$ws_CRP.Range("N36").Value = -6775
# row 40
$ws_CRP.Range("H40").Value = 3024.5
$ws_CRP.Range("I40").Value = 50
$ws_CRP.Range("J40").Value = 5999
$ws_CRP.Range("K40").Value = 50
$ws_CRP.Range("L40").Value = 5999
$ws_CRP.Range("M40").Value = 110
$ws_CRP.Range("N40").Value = -6319
# row 132
$ws_CRP.Range("H132").Value = 3912
$ws_CRP.Range("I132").Value = 891.6
$ws_CRP.Range("K132").Value = 2674.8
$ws_CRP.Range("M132").Value = -144.8000000000002
# row 133
$ws_CRP.Range("H133").Value = 59996.25
$ws_CRP.Range("I133").Value = 0
$ws_CRP.Range("K133").Value = 0
$ws_CRP.Range("M133").ClearContents()
# row 134
$ws_CRP.Range("H134").Value = 2539.0334
$ws_CRP.Range("I134").Value = 1973.1072
$ws_CRP.Range("J134").Value = 10462
$ws_CRP.Range("K134").Value = 5919.321599999999
$ws_CRP.Range("L134").Value = 31386
$ws_CRP.Range("M134").Value = -3384.321599999999
$ws_CRP.Range("N134").Value = -36456

$ws_CUL = $wb.Worksheets.Item("CUL")
# row 58
$ws_CUL.Range("H58").Value = 2126
$ws_CUL.Range("I58").Value = 2126
$ws_CUL.Range("K58").Value = 6378
$ws_CUL.Range("M58").Value = -6250
# row 69
$ws_CUL.Range("H69").Value = 6627.3335
$ws_CUL.Range("I69").Value = 777
$ws_CUL.Range("K69").Value = 2331
$ws_CUL.Range("M69").Value = -1520
# row 72
$ws_CUL.Range("H72").Value = 6627.3335
$ws_CUL.Range("I72").Value = 777
$ws_CUL.Range("K72").Value = 6993
$ws_CUL.Range("M72").Value = -2937
# row 93
$ws_CUL.Range("H93").Value = 6800
$ws_CUL.Range("J93").Value = 6800
$ws_CUL.Range("L93").Value = 20400
$ws_CUL.Range("N93").Value = -24144
# row 114
$ws_CUL.Range("H114").Value = 1223
$ws_CUL.Range("I114").Value = 1083
$ws_CUL.Range("J114").Value = 1293
$ws_CUL.Range("K114").Value = 3249
$ws_CUL.Range("L114").Value = 3879
$ws_CUL.Range("M114").Value = 5
$ws_CUL.Range("N114").Value = -10387
# row 128
$ws_CUL.Range("H128").Value = 1875988
$ws_CUL.Range("I128").Value = 1875988
$ws_CUL.Range("K128").Value = 5627964
$ws_CUL.Range("M128").Value = -5622984
# row 132
$ws_CUL.Range("H132").Value = 3382.25
$ws_CUL.Range("J132").Value = 5777.6665
$ws_CUL.Range("L132").Value = 51998.9985
$ws_CUL.Range("N132").Value = -57058.9985

$ws_GSM = $wb.Worksheets.Item("GSM")
# row 24
$ws_GSM.Range("H24").Value = 13749.75
$ws_GSM.Range("I24").Value = 8333.333000000001
$ws_GSM.Range("J24").Value = 29999
$ws_GSM.Range("K24").Value = 8333.333000000001
$ws_GSM.Range("L24").Value = 29999
$ws_GSM.Range("M24").Value = -8160.333000000001
$ws_GSM.Range("N24").Value = -30345
# row 102
$ws_GSM.Range("H102").Value = 9529241
$ws_GSM.Range("I102").Value = 13336938
$ws_GSM.Range("K102").Value = 13336938
$ws_GSM.Range("M102").Value = -13335316
# row 132
$ws_GSM.Range("H132").Value = 1314908.8
$ws_GSM.Range("I132").Value = 1640882.1
$ws_GSM.Range("J132").Value = 11015.5
$ws_GSM.Range("K132").Value = 4922646.300000001
$ws_GSM.Range("L132").Value = 33046.5
$ws_GSM.Range("M132").Value = -4920116.300000001
$ws_GSM.Range("N132").Value = -38106.5

$ws_LTW = $wb.Worksheets.Item("LTW")
# row 46
$ws_LTW.Range("H46").Value = 2467.4119
$ws_LTW.Range("I46").Value = 600.25
$ws_LTW.Range("J46").Value = 3041.923
$ws_LTW.Range("K46").Value = 600.25
$ws_LTW.Range("L46").Value = 3041.923
$ws_LTW.Range("M46").Value = -412.25
$ws_LTW.Range("N46").Value = -3417.923
# row 97
$ws_LTW.Range("H97").Value = 8545.454
$ws_LTW.Range("I97").Value = 7000
$ws_LTW.Range("J97").Value = 8619.048000000001
$ws_LTW.Range("K97").Value = 7000
$ws_LTW.Range("L97").Value = 8619.048000000001
$ws_LTW.Range("M97").Value = -6009
$ws_LTW.Range("N97").Value = -10601.048
# row 132
$ws_LTW.Range("H132").Value = 6462.4688
$ws_LTW.Range("I132").Value = 3339.5
$ws_LTW.Range("K132").Value = 10018.5
$ws_LTW.Range("M132").Value = -7488.5

$ws_WVR = $wb.Worksheets.Item("WVR")
# row 4
$ws_WVR.Range("H4").Value = 10000
$ws_WVR.Range("J4").Value = 10000
$ws_WVR.Range("L4").Value = 10000
$ws_WVR.Range("N4").Value = -10226
